$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Output"

# New data rows: year, gender (M/F), value
$data = @(
    @(1980, "M", 5.6),
    @(1980, "F", 3.7),
    @(1983, "M", 6.4),
    @(1983, "F", 4.7),
    @(1984, "M", 6.6),
    @(1984, "F", 5),
    @(1985, "M", 6.5),
    @(1985, "F", 4.9000000000000004),
    @(1986, "M", 6.6),
    @(1986, "F", 5),
    @(1987, "M", 6.6),
    @(1987, "F", 5.2),
    @(1988, "M", 6.9),
    @(1988, "F", 5.5),
    @(1989, "M", 7),
    @(1989, "F", 5.6),
    @(1990, "M", 7.3),
    @(1990, "F", 5.9),
    @(1991, "M", 7.6),
    @(1991, "F", 6.2),
    @(1992, "M", 7.8),
    @(1992, "F", 6.3),
    @(1993, "M", 7.9),
    @(1993, "F", 6.5),
    @(1994, "M", 8.1),
    @(1994, "F", 6.6),
    @(1995, "M", 8.4),
    @(1995, "F", 7.2),
    @(1996, "M", 8.4),
    @(1996, "F", 7.1),
    @(1997, "M", 8.6),
    @(1997, "F", 7.2),
    @(1998, "M", 8.8000000000000007),
    @(1998, "F", 7.4),
    @(1999, "M", 8.8000000000000007),
    @(1999, "F", 7.6),
    @(2000, "M", 9.1999999999999993),
    @(2000, "F", 8.1),
    @(2001, "M", 9.1999999999999993),
    @(2001, "F", 7.9),
    @(2002, "M", 9.3000000000000007),
    @(2002, "F", 8.1),
    @(2003, "M", 9.4),
    @(2003, "F", 8.1999999999999993),
    @(2004, "M", 9.6),
    @(2004, "F", 8.3000000000000007),
    @(2005, "M", 9.9),
    @(2005, "F", 8.8000000000000007),
    @(2006, "M", 9.9),
    @(2006, "F", 8.6999999999999993),
    @(2007, "M", 9.9),
    @(2007, "F", 8.8000000000000007),
    @(2008, "M", 10.199999999999999),
    @(2008, "F", 9.1),
    @(2009, "M", 10.3),
    @(2009, "F", 9.1999999999999993),
    @(2010, "M", 10.6),
    @(2010, "F", 9.6999999999999993),
    @(2011, "M", 10.7),
    @(2011, "F", 9.6999999999999993),
    @(2012, "M", 10.9),
    @(2012, "F", 9.9),
    @(2013, "M", 11),
    @(2013, "F", 10),
    @(2014, "M", 11.1),
    @(2014, "F", 10.1),
    @(2015, "M", 11.2),
    @(2015, "F", 10.3),
    @(2016, "M", 11.2),
    @(2016, "F", 10.3),
    @(2017, "M", 11.3),
    @(2017, "F", 10.4),
    @(2018, "M", 11.6),
    @(2018, "F", 10.6)
)

# Write header row
$ws.Cells.Item(1, 1).Value = "year"
$ws.Cells.Item(1, 2).Value = "gender"
$ws.Cells.Item(1, 3).Value = "value"

# Write data rows starting at row 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-unused column D (old "value" column)
$ws.Columns.Item(4).Delete()
